$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")
$ws.Activate()

# Duplicate the formatting of the row above (s="4"/s="5") onto the new row
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# Fill in the new data row: type "string", name "is_override"
$ws.Range("A12").Value = "string"
$ws.Range("B12").Value = "is_override"

$ws.Range("B12").Select()
